$wb = $excel.ActiveWorkbook
$wsGermany = $wb.Worksheets.Item("Germany")
$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsCzech   = $wb.Worksheets.Item("Czech")

# --- Update test-data markers (B4) for Germany and Czech markets ---
# Order matters for how the shared-string table gets rebuilt on save,
# so update Belgium before Germany (matches the authored commit).
$wsBelgium.Range("B4").Value = "NGC-3478/T2265/T2264/T2267/T2268"
$wsGermany.Range("B4").Value = "NGC-3475/T1730/T1746"

# --- Germany: remove the "PROFILE Communicator" and "RS800" printer rows ---
$wsGermany.Rows("12").Delete()
$wsGermany.Rows("13").Delete()

# --- Belgium: remove the "RS800" printer row (adds Belgium market test data) ---
$wsBelgium.Rows("13").Delete()

# --- Selection / active sheet bookkeeping: Germany becomes the active tab ---
$wsBelgium.Range("A11").Select()
$wsGermany.Range("A11").Select()
